$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 24

$ws.Cells.Item($row, 1).Value = "VB3KB3"
$ws.Cells.Item($row, 2).Value = "Cinta Flex Scanner Samsung"
$ws.Cells.Item($row, 3).Value = "CLX3300 CLX3305 M2070 M2876 M3370 M3375 M3870 M3875 M4070 M4075 SCX3400 SCX3405 SCX4833 SCX4835"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 100000
$ws.Cells.Item($row, 6).Value = 3
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E24-D24)*G24"
$ws.Cells.Item($row, 9).Formula = "=D24*F24"
$ws.Cells.Item($row, 10).Value = 0
